$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6, column A: switch its number format from date-only ("YYYY-MM-DD")
# to the date-time format used by the other data rows ("YYYY-MM-DD HH:MM:SS").
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 7 - mirrors the layout/format of the existing data rows, with A7
# using the date-only format (like the previous "last" row, A6, used to).
$ws.Range("A7").Value = 45856
$ws.Range("A7").NumberFormat = "YYYY-MM-DD"

$ws.Range("B7").Value = "ert5678"
$ws.Range("C7").Value = "Diego"
$ws.Range("D7").Value = "2025-07-18 13:42:59"
$ws.Range("E7").Value = "2025-07-18 13:42:59"
$ws.Range("F7").Value = "2025-07-18 13:43:01"
$ws.Range("G7").Value = "2025-07-18 13:43:03"
$ws.Range("H7").Value = "2025-07-18 13:43:04"
$ws.Range("I7").Value = "2025-07-18 13:43:05"
$ws.Range("J7").Value = "2025-07-18 13:43:06"
$ws.Range("K7").Value = "0:00:02"
$ws.Range("L7").Value = "0:00:00"
$ws.Range("M7").Value = "0:00:07"
# N7 stays blank (mirrors N2:N6, which are empty inline strings).
$ws.Range("O7").Value = "2025-07-18 13:43:08"
$ws.Range("P7").Value = "2025-07-18 13:43:09"
$ws.Range("Q7").Value = "2025-07-18 13:43:10"
$ws.Range("R7").Value = "2025-07-18 13:43:11"
$ws.Range("S7").Value = "0:00:01"
$ws.Range("T7").Value = "0:00:01"
$ws.Range("U7").Value = "0:00:04"
$ws.Range("V7").Value = "0:00:01"
$ws.Range("W7").Value = "2025-07-18 13:43:07"
